$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.208.69'
$ws.Range('E2').Value = '  +1.39%  '
$ws.Range('D3').Value = '1.655.40'
$ws.Range('E3').Value = '  +0.62%  '
$ws.Range('E4').Value = '  -0.62%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '220.01'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.58%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.502'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.14%  '
$ws.Range('E7').Value = '  -0.66%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.254'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.76%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.0626'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.32%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.60'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.99%  '
$ws.Range('E11').Value = '  +0.47%  '
$ws.Range('D12').Value = '1.884.74'
$ws.Range('E12').Value = '  +0.50%  '
$ws.Range('D13').Value = '1.647.61'
$ws.Range('E13').Value = '  -0.20%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.20'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.91%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.531'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.31%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '65.98'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +1.90%  '
$ws.Range('D17').Value = '27.198.44'
$ws.Range('E17').Value = '  +1.37%  '
$ws.Range('D18').Value = '0.0₃0738'
$ws.Range('E18').Value = '  +0.34%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '221.08'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +3.08%  '
$ws.Range('E20').Value = '  -0.72%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.73'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +7.57%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.43'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.82%  '
$ws.Range('E23').Value = '  -1.42%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '9.28'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.78%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '147.59'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.36%  '
$ws.Range('E26').Value = '  -0.59%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.35'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +2.17%  '
$ws.Range('E28').Value = '  +0.19%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '15.97'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +1.94%  '
$ws.Range('E30').Value = '  +1.45%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.20'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +1.15%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.37'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.34%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.01'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.30%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.57'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +2.62%  '
$ws.Range('D35').Value = '1.266.93'
$ws.Range('E35').Value = '  -2.20%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.45'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.01%  '
$ws.Range('E37').Value = '  -1.58%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.539'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.83%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.827'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.58%  '
$ws.Range('E40').Value = '  -0.72%  '
$ws.Range('E41').Value = '  +0.00%  '
$ws.Range('E42').Value = '  +1.04%  '
$ws.Range('D43').Value = '1.793.59'
$ws.Range('E43').Value = '  +0.42%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '61.89'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.58%  '
$ws.Range('B45').Value = 'MXToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.09'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -6.32%  '
$ws.Range('B46').Value = 'Quant'
$ws.Range('C46').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '92.62'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.79%  '
$ws.Range('E47').Value = '  +0.48%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0518'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.57%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '7.68'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.01%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0977'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.63%  '
$ws.Range('E51').Value = '  -0.08%  '
